# RPA datasets push 2023-12-23
# Insert a new IPO record row (IBKS제23호스팩 / IBK) right after the header+first
# data row, i.e. as the new row 3, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 3..29 down to 4..30, leaving a blank row 3 behind.
$ws.Rows.Item(3).Insert()

# Populate the new row 3 with the new offering's data.
$ws.Range("A3").Value = "2023-12-04"
$ws.Range("B3").Value = "2023-12-08"
$ws.Range("C3").Value = "2023-12-22"
$ws.Range("D3").Value = "IBK"
$ws.Range("E3").Value = "IBKS제23호스팩"
$ws.Range("F3").Value = 4000000
$ws.Range("G3").Value = 4000000
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 4230000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 2000
$ws.Range("N3").Value = "855.86"
$ws.Range("O3").Value = "-"
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0
$ws.Range("W3").Value = 0
$ws.Range("X3").Value = 0
$ws.Range("Y3").Value = "기업인수 및 합병"
